# open Dingo dataprovider parallel
#
# Add a new batch_016 test case as row 17 of the batchsql cases sheet,
# matching the existing rows' layout:
#   A=TestID  B=Testable  C=Title  D=Component  E=Sub_component
#   G=Table_name  I=Batch_sql  J=Query_sql1  K=Query_result1  N=Validation_type

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row by duplicating the previous (last) data row so that the
# new row picks up the same row-level formatting/styles used throughout the
# sheet, then overwrite the cells that actually change for batch_016.
$ws.Rows("16").Copy()
$ws.Rows("17").Insert(-4121)  # xlShiftDown

$ws.Range("A17").Value = "batch_016"
$ws.Range("C17").Value = "批量操作语句16执行"
$ws.Range("D17").Value = "batchsql"
$ws.Range("G17").Value = "v"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("I17").Value = "batch_sql_016"
$ws.Range("J17").Value = "select a.id from v as a "
$ws.Range("K17").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_016.csv"
$ws.Range("N17").NumberFormat = "@"

# B17, E17 and N17 already carry the right text ("y", "SingleTable",
# "csv_containsAll") copied over from row 16, so they don't need to be
# retyped.

$ws.Range("J17").Select()
